$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.050374230381654
$ws.Range("D2").Value = 1.061073235376192
$ws.Range("E2").Value = 1.047415766775655
$ws.Range("F2").Value = 1.067764813040254
$ws.Range("I2").Value = 1.046890175727271
$ws.Range("J2").Value = 1.055407722577377
$ws.Range("K2").Value = 1.063798287502832
$ws.Range("L2").Value = 1.050178466896383
$ws.Range("M2").Value = 1.070471799631482
$ws.Range("N2").Value = 1.056906522600394

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.051758728536194
$ws.Range("D3").Value = 1.062001342463486
$ws.Range("E3").Value = 1.048609969482911
$ws.Range("F3").Value = 1.068940439373152
$ws.Range("I3").Value = 1.047278141525119
$ws.Range("J3").Value = 1.056439752721359
$ws.Range("K3").Value = 1.064540657860557
$ws.Range("L3").Value = 1.05118351795605
$ws.Range("M3").Value = 1.071462379341765
$ws.Range("N3").Value = 1.057940018345557

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.052653807882161
$ws.Range("D4").Value = 1.062601358422196
$ws.Range("E4").Value = 1.049382206800262
$ws.Range("F4").Value = 1.069700937488646
$ws.Range("I4").Value = 1.047527620201299
$ws.Range("J4").Value = 1.057106310577251
$ws.Range("K4").Value = 1.06501986147795
$ws.Range("L4").Value = 1.051832793764083
$ws.Range("M4").Value = 1.072102560726352
$ws.Range("N4").Value = 1.05860752279006

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.053029916246423
$ws.Range("D5").Value = 1.062853479300421
$ws.Range("E5").Value = 1.049706740777669
$ws.Range("F5").Value = 1.070020602654865
$ws.Range("I5").Value = 1.047632128383797
$ws.Range("J5").Value = 1.057386239015731
$ws.Range("K5").Value = 1.065221042517914
$ws.Range("L5").Value = 1.052105498548225
$ws.Range("M5").Value = 1.072371505916007
$ws.Range("N5").Value = 1.05888784875903

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.053093055898058
$ws.Range("D6").Value = 1.062895804166542
$ws.Range("E6").Value = 1.049761224831315
$ws.Range("F6").Value = 1.070074273077971
$ws.Range("I6").Value = 1.047649653946441
$ws.Range("J6").Value = 1.05743322317064
$ws.Range("K6").Value = 1.06525480554786
$ws.Range("L6").Value = 1.052151272259498
$ws.Range("M6").Value = 1.072416652086493
$ws.Range("N6").Value = 1.058934899636828

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.052658834173863
$ws.Range("D7").Value = 1.062604727766426
$ws.Range("E7").Value = 1.04938654368554
$ws.Range("F7").Value = 1.069705209057054
$ws.Range("I7").Value = 1.047529018107672
$ws.Range("J7").Value = 1.057110052140169
$ws.Range("K7").Value = 1.065022550752166
$ws.Range("L7").Value = 1.051836438640539
$ws.Range("M7").Value = 1.072106155117637
$ws.Range("N7").Value = 1.058611269666426

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.050842293017347
$ws.Range("D8").Value = 1.061387004375528
$ws.Range("E8").Value = 1.047819456860727
$ws.Range("F8").Value = 1.068162166430128
$ws.Range("I8").Value = 1.047021614484425
$ws.Range("J8").Value = 1.055756759506768
$ws.Range("K8").Value = 1.064049416034929
$ws.Range("L8").Value = 1.050518349697924
$ws.Range("M8").Value = 1.070806735476143
$ws.Range("N8").Value = 1.057256055202261

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.047635108630005
$ws.Range("D9").Value = 1.059237087622264
$ws.Range("E9").Value = 1.045054143650787
$ws.Range("F9").Value = 1.065441414623123
$ws.Range("I9").Value = 1.046115499042419
$ws.Range("J9").Value = 1.053362489897447
$ws.Range("K9").Value = 1.062325680829182
$ws.Range("L9").Value = 1.048187465801651
$ws.Range("M9").Value = 1.068510849922197
$ws.Range("N9").Value = 1.054858385455468

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.045492523133558
$ws.Range("D10").Value = 1.057800940769615
$ws.Range("E10").Value = 1.043207777050172
$ws.Range("F10").Value = 1.063626286033895
$ws.Range("I10").Value = 1.045503280526654
$ws.Range("J10").Value = 1.051759662552449
$ws.Range("K10").Value = 1.061170413710465
$ws.Range("L10").Value = 1.046627819218497
$ws.Range("M10").Value = 1.066975998506249
$ws.Range("N10").Value = 1.053253281911793

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.044563636547135
$ws.Range("D11").Value = 1.05717837183402
$ws.Range("E11").Value = 1.042407565210434
$ws.Range("F11").Value = 1.062839975953862
$ws.Range("I11").Value = 1.045236236137125
$ws.Range("J11").Value = 1.051063999581259
$ws.Range("K11").Value = 1.060668699276563
$ws.Range("L11").Value = 1.045951076266761
$ws.Range("M11").Value = 1.066310350221095
$ws.Range("N11").Value = 1.05255663101939

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.044218430823565
$ws.Range("D12").Value = 1.056947013751763
$ws.Range("E12").Value = 1.042110218347023
$ws.Range("F12").Value = 1.062547850289398
$ws.Range("I12").Value = 1.04513674959626
$ws.Range("J12").Value = 1.050805350861749
$ws.Range("K12").Value = 1.060482116458609
$ws.Range("L12").Value = 1.045699488986054
$ws.Range("M12").Value = 1.066062938752827
$ws.Range("N12").Value = 1.052297614989031

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.04429248667379
$ws.Range("D13").Value = 1.056996645775589
$ws.Range("E13").Value = 1.042174005415621
$ws.Range("F13").Value = 1.062610514775687
$ws.Range("I13").Value = 1.045158103140796
$ws.Range("J13").Value = 1.050860843198968
$ws.Range("K13").Value = 1.060522149262741
$ws.Range("L13").Value = 1.04575346509115
$ws.Range("M13").Value = 1.066116016638817
$ws.Range("N13").Value = 1.052353186131734

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.044535105355307
$ws.Range("D14").Value = 1.057159249922338
$ws.Range("E14").Value = 1.042382988744974
$ws.Range("F14").Value = 1.062815829895939
$ws.Range("I14").Value = 1.04522801856192
$ws.Range("J14").Value = 1.051042624689768
$ws.Range("K14").Value = 1.060653280858229
$ws.Range("L14").Value = 1.045930284386897
$ws.Range("M14").Value = 1.066289902392647
$ws.Range("N14").Value = 1.052535225773102

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.044684567255428
$ws.Range("D15").Value = 1.057259421300896
$ws.Range("E15").Value = 1.042511735238223
$ws.Range("F15").Value = 1.062942323930897
$ws.Range("I15").Value = 1.04527105671104
$ws.Range("J15").Value = 1.051154593217604
$ws.Range("K15").Value = 1.060734045655613
$ws.Range("L15").Value = 1.046039200016663
$ws.Range("M15").Value = 1.066397017854852
$ws.Range("N15").Value = 1.052647353309089

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.045554145923651
$ws.Range("D16").Value = 1.057842243542784
$ws.Range("E16").Value = 1.043260868907142
$ws.Range("F16").Value = 1.063678463277866
$ws.Range("I16").Value = 1.045520962182018
$ws.Range("J16").Value = 1.051805796776488
$ws.Range("K16").Value = 1.061203679577034
$ws.Range("L16").Value = 1.046672702493968
$ws.Range("M16").Value = 1.06702015309256
$ws.Range("N16").Value = 1.053299481651722

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.046099302331463
$ws.Range("D17").Value = 1.058207641675963
$ws.Range("E17").Value = 1.043730584340429
$ws.Range("F17").Value = 1.064140128988907
$ws.Range("I17").Value = 1.045677198269091
$ws.Range("J17").Value = 1.05221384090994
$ws.Range("K17").Value = 1.061497871939126
$ws.Range("L17").Value = 1.047069702962909
$ws.Range("M17").Value = 1.067410747031028
$ws.Range("N17").Value = 1.053708105254646

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.046417173915474
$ws.Range("D18").Value = 1.058420704056576
$ws.Range("E18").Value = 1.044004491923297
$ws.Range("F18").Value = 1.064409377463675
$ws.Range("I18").Value = 1.045768140112615
$ws.Range("J18").Value = 1.052451689321922
$ws.Range("K18").Value = 1.061669327073935
$ws.Range("L18").Value = 1.047301131237981
$ws.Range("M18").Value = 1.067638472982176
$ws.Range("N18").Value = 1.053946291438655

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.046525541606415
$ws.Range("D19").Value = 1.058493341316927
$ws.Range("E19").Value = 1.044097875663897
$ws.Range("F19").Value = 1.064501178670109
$ws.Range("I19").Value = 1.045799117082406
$ws.Range("J19").Value = 1.052532763000692
$ws.Range("K19").Value = 1.061727764767672
$ws.Range("L19").Value = 1.047380019445144
$ws.Range("M19").Value = 1.067716104615249
$ws.Range("N19").Value = 1.054027480251348

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.04604082351754
$ws.Range("D20").Value = 1.05816844498009
$ws.Range("E20").Value = 1.043680195536513
$ws.Range("F20").Value = 1.064090600096205
$ws.Range("I20").Value = 1.045660455074593
$ws.Range("J20").Value = 1.052170077887001
$ws.Range("K20").Value = 1.061466322615926
$ws.Range("L20").Value = 1.04702712260035
$ws.Range("M20").Value = 1.067368850453435
$ws.Range("N20").Value = 1.053664280083195

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.044463665101776
$ws.Range("D21").Value = 1.057111370057839
$ws.Range("E21").Value = 1.042321451502173
$ws.Range("F21").Value = 1.062755371254466
$ws.Range("I21").Value = 1.045207438350493
$ws.Range("J21").Value = 1.050989101463288
$ws.Range("K21").Value = 1.060614672063841
$ws.Range("L21").Value = 1.045878221470622
$ws.Range("M21").Value = 1.066238701812128
$ws.Range("N21").Value = 1.0524816265375

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.043471022551425
$ws.Range("D22").Value = 1.056446117784541
$ws.Range("E22").Value = 1.041466502011022
$ws.Range("F22").Value = 1.061915539704787
$ws.Range("I22").Value = 1.044920904935527
$ws.Range("J22").Value = 1.050245134917527
$ws.Range("K22").Value = 1.060077910096316
$ws.Range("L22").Value = 1.045154617058369
$ws.Range("M22").Value = 1.065527205123488
$ws.Range("N22").Value = 1.05173660347391

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.043997339778032
$ws.Range("D23").Value = 1.056798840705002
$ws.Range("E23").Value = 1.041919790149747
$ws.Range("F23").Value = 1.062360781619238
$ws.Range("I23").Value = 1.045072963704943
$ws.Range("J23").Value = 1.050639663444488
$ws.Range("K23").Value = 1.060362581224843
$ws.Range("L23").Value = 1.045538332644497
$ws.Range("M23").Value = 1.065904471784832
$ws.Range("N23").Value = 1.052131692276633

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.046067247912584
$ws.Range("D24").Value = 1.058186156490395
$ws.Range("E24").Value = 1.043702964282554
$ws.Range("F24").Value = 1.064112980173785
$ws.Range("I24").Value = 1.045668021184396
$ws.Range("J24").Value = 1.052189852997679
$ws.Range("K24").Value = 1.061480578836685
$ws.Range("L24").Value = 1.04704636325147
$ws.Range("M24").Value = 1.067387782025087
$ws.Range("N24").Value = 1.053684083276799

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.048465007932631
$ws.Range("D25").Value = 1.05979339163353
$ws.Range("E25").Value = 1.045769526631682
$ws.Range("F25").Value = 1.066145011804363
$ws.Range("I25").Value = 1.046351181468193
$ws.Range("J25").Value = 1.05398262320949
$ws.Range("K25").Value = 1.062772377710316
$ws.Range("L25").Value = 1.048791050651567
$ws.Range("M25").Value = 1.069105132721935
$ws.Range("N25").Value = 1.055479399427945
